$d = $word.ActiveDocument

# 1) "Task 3: Using the platform..." -> "Task 2.3: Using the platform..."
#    (the objectives list item becomes "Task 2.3" to match the later numbering)
$r1 = $d.Content
$r1.Find.Execute("Task 3: Using the platform", $true, $false, $false, $false, $false, $true, 1, $false, "Task 2.3: Using the platform", 2)

# 2) "d =" -> "d (private key) ="  (annotate the private key exponent d).
#    This exact three-character paragraph ("d =") only occurs once (paragraph 73);
#    a document-wide Find would also hit the unrelated "(d = c^d mod n)" phrase
#    later in the report, so scope the Find to that specific paragraph.
$p2 = $d.Paragraphs.Item(73)
$p2.Range.Find.Execute("d =", $true, $false, $false, $false, $false, $true, 1, $false, "d (private key) =", 2)

# 3) "The public key, h, is computed using the formula:" -> "The public key, is computed using the formula:"
$r3 = $d.Content
$r3.Find.Execute("The public key, h, is computed using the formula:", $true, $false, $false, $false, $false, $true, 1, $false, "The public key, is computed using the formula:", 2)

# 4) "Thus both actors calculate the shared key:" -> "Thus both actors calculate the shared key (KpubB^KprA mod p):"
$r4 = $d.Content
$r4.Find.Execute("Thus both actors calculate the shared key:", $true, $false, $false, $false, $false, $true, 1, $false, "Thus both actors calculate the shared key (KpubB^KprA mod p):", 2)
